# Applies the edit described by the diff:
#  - exp_10 (sheet10): add "UDP" tech label in col A for rows 2-10, and
#    add min/mean/max conn ms data in cols H/I/J for rows 2-10.
#  - Update selection/active-cell bookkeeping on exp_1, exp_7, exp_9, exp_10
#    and make exp_1 the active (selected) sheet tab.

$wb = $excel.ActiveWorkbook

$ws1  = $wb.Worksheets.Item("exp_1")
$ws7  = $wb.Worksheets.Item("exp_7")
$ws9  = $wb.Worksheets.Item("exp_9")
$ws10 = $wb.Worksheets.Item("exp_10")

# --- exp_10: fill in new "UDP" tech column for rows 2-10 ---
$ws10.Range("A2").Value = "UDP"
$ws10.Range("A3").Value = "UDP"
$ws10.Range("A4").Value = "UDP"
$ws10.Range("A5").Value = "UDP"
$ws10.Range("A6").Value = "UDP"
$ws10.Range("A7").Value = "UDP"
$ws10.Range("A8").Value = "UDP"
$ws10.Range("A9").Value = "UDP"
$ws10.Range("A10").Value = "UDP"

# --- exp_10: fill in new min/mean/max conn ms columns (H/I/J) ---
$ws10.Range("H2").Value = 0.030956029891967701
$ws10.Range("I2").Value = 0.031160253744858899
$ws10.Range("J2").Value = 0.032351970672607401

$ws10.Range("H3").Value = 0.030949592590332
$ws10.Range("I3").Value = 0.031083524227142299
$ws10.Range("J3").Value = 0.031283617019653299

$ws10.Range("H4").Value = 0.030978918075561499
$ws10.Range("I4").Value = 0.031113496193519
$ws10.Range("J4").Value = 0.031389951705932603

$ws10.Range("H5").Value = 0.0308954715728759
$ws10.Range("I5").Value = 0.031074872383704499
$ws10.Range("J5").Value = 0.031439781188964802

$ws10.Range("H6").Value = 0.030788183212280201
$ws10.Range("I6").Value = 0.030977189540862999
$ws10.Range("J6").Value = 0.031207561492919901

$ws10.Range("H7").Value = 0.030765533447265601
$ws10.Range("I7").Value = 0.030928405431600699
$ws10.Range("J7").Value = 0.031269073486328097

$ws10.Range("H8").Value = 0.030736684799194301
$ws10.Range("I8").Value = 0.030900579232435901
$ws10.Range("J8").Value = 0.031119823455810498

$ws10.Range("H9").Value = 0.030732154846191399
$ws10.Range("I9").Value = 0.030891730235173099
$ws10.Range("J9").Value = 0.031092405319213801

$ws10.Range("H10").Value = 0.030709505081176699
$ws10.Range("I10").Value = 0.0308484205832848
$ws10.Range("J10").Value = 0.031126022338867101

# --- update selections on the sheets that were visited, in diff order ---
$ws10.Activate()
$ws10.Range("H15").Select() | Out-Null

$ws7.Activate()
$ws7.Range("H2").Select() | Out-Null

$ws9.Activate()
$ws9.Range("I17").Select() | Out-Null

$ws1.Activate()
$ws1.Range("E29").Select() | Out-Null
